$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 158, shifting existing rows 158:226 down to 159:227
$ws.Rows.Item(158).Insert()

# Populate the newly inserted row 158 with the new record's data
$ws.Range("A158").Value = 10
$ws.Range("B158").Value = "Vega Modelo de Temuco"
$ws.Range("C158").Value = "La Araucanía"
$ws.Range("D158").Value = 44784
$ws.Range("E158").Value = 9
$ws.Range("F158").Value = 100112005
$ws.Range("G158").Value = "Puerro"
$ws.Range("H158").Value = "Azul de Maquehue"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 80
$ws.Range("K158").Value = 15000
$ws.Range("L158").Value = 16000
$ws.Range("M158").Value = 15500
$ws.Range("N158").Value = "`$/docena de paquetes"
$ws.Range("O158").Value = "Provincia de Cautín"
$ws.Range("P158").Value = 1292
$ws.Range("Q158").Value = 12
$ws.Range("R158").Value = "Hortaliza"
